# Applies the cryptos-list price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.867.65"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "1.624.97"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'" + "210.98"
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").Value = "'" + "23.45"
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("E9").Value = "  -2.14%  "
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").Value = "1.856.18"
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("D13").Value = "1.624.18"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").Value = "'" + "4.02"
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("D15").Value = "'" + "0.561"
$ws.Range("E15").Value = "  -2.35%  "
$ws.Range("D16").Value = "'" + "65.37"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").Value = "27.856.47"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "'" + "229.40"
$ws.Range("E18").Value = "  -1.37%  "
$ws.Range("D19").Value = "'" + "7.65"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("E22").Value = "  -0.93%  "
$ws.Range("E23").Value = "  -5.54%  "
$ws.Range("E24").Value = "  -2.59%  "
$ws.Range("D25").Value = "'" + "154.72"
$ws.Range("E25").Value = "  +2.07%  "
$ws.Range("E26").Value = "  -1.35%  "
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("D28").Value = "'" + "15.52"
$ws.Range("E28").Value = "  -1.34%  "
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("E32").Value = "  +1.88%  "
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("D34").Value = "1.399.92"
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").Value = "'" + "0.999"
$ws.Range("E36").Value = "  +9.28%  "
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("D40").Value = "'" + "0.860"
$ws.Range("E40").Value = "  -2.92%  "
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("E43").Value = "  -2.07%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "'" + "65.92"
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("D47").Value = "1.766.15"
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("D48").Value = "'" + "87.97"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("D49").Value = "'" + "0.103"
$ws.Range("E49").Value = "  +2.19%  "
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("D51").Value = "'" + "7.54"
$ws.Range("E51").Value = "  -1.24%  "
